# Auto-generated script applying the 2023-03-11 daily update
# to the "violent-crime-full-year" workbook (column J = 2023 YTD totals).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1210
$ws.Range("J3").Value = 1290
$ws.Range("J4").Value = 281
$ws.Range("J5").Value = 94
$ws.Range("J6").Value = 1701
$ws.Range("J7").Value = 4576

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 52
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 70
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 30
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 14
$ws.Range("J7").Value = 130
$ws.Range("J8").Value = 288
$ws.Range("J12").Value = 15
$ws.Range("J15").Value = 51
$ws.Range("J17").Value = 12
$ws.Range("J19").Value = 161
$ws.Range("J20").Value = 97
$ws.Range("J29").Value = 255
$ws.Range("J33").Value = 189
$ws.Range("J37").Value = 157
$ws.Range("J42").Value = 188
$ws.Range("J44").Value = 37
$ws.Range("J46").Value = 17
$ws.Range("J48").Value = 30
$ws.Range("J50").Value = 23
$ws.Range("J52").Value = 101
$ws.Range("J53").Value = 47
$ws.Range("J54").Value = 91
$ws.Range("J55").Value = 58
$ws.Range("J57").Value = 18
$ws.Range("I63").Value = 186
$ws.Range("J63").Value = 18
$ws.Range("J65").Value = 115
$ws.Range("J66").Value = 10
$ws.Range("J67").Value = 170
$ws.Range("J71").Value = 21
$ws.Range("J76").Value = 76
$ws.Range("J78").Value = 61
$ws.Range("J79").Value = 135
$ws.Range("J82").Value = 8
$ws.Range("J83").Value = 106
$ws.Range("J85").Value = 193
$ws.Range("J88").Value = 35
$ws.Range("J89").Value = 53
$ws.Range("J90").Value = 50
$ws.Range("I91").Value = 279
$ws.Range("J91").Value = 64
$ws.Range("J93").Value = 20
$ws.Range("J94").Value = 35
$ws.Range("J97").Value = 27
$ws.Range("J101").Value = 4576

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 33
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 46
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 14
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 96
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 40
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 52
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 34
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 17
$ws.Range("J3").Value = 27
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 279
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 24
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 7
$ws.Range("J4").Value = 2
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J3").Value = 1
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 97
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J3").Value = 12
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 8

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 15
